$p = $ppt.ActivePresentation
Write-Output $p.NotesMaster
Write-Output $p.NotesMaster.Theme
